# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@79dae63976c2f2fd3ec2efcb4d3966d843ed4800
# This script bumps the IG version/date metadata on the "Metadata" sheet and
# removes the now-obsolete "Extension.extension.extension.*" slice rows from
# the "Elements" (StructureDefinition) table, since the extension definition
# was simplified to drop a level of nesting.

$wb = $excel.ActiveWorkbook

# --- Sheet 1 ("Metadata"): bump Version and Date -----------------------
$wsMeta = $wb.Worksheets.Item(1)
$wsMeta.Range("B3").Value2 = "6.1.0"
$wsMeta.Range("B8").Value2 = "2022-05-31T20:10:14+00:00"

# --- Sheet 2 ("Elements"): drop the Extension.extension.extension.* rows
$wsElements = $wb.Worksheets.Item(2)

# Rows 9-12 held the nested Extension.extension.extension.id/.extension/.url/
# .value[x] definitions. That nesting level no longer exists, so delete the
# whole block; the rows below (old 13-17, Extension.extension.url /
# Extension.extension.value[x] / Extension.url / Extension.value[x]) shift
# up to become the new rows 9-13.
$wsElements.Range("A9:A12").EntireRow.Delete()

# The first column's "best fit" width shrinks now that the longest path in
# it is shorter than before.
$wsElements.Columns.Item(1).AutoFit()
